$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F1 - copy the style of E1 (existing header style) and set the label
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

# time_taken values for each data row, matching the diff
$times = @(
    "2021-10-05 13:40:35.393994",
    "2021-10-05 13:40:35.394006",
    "2021-10-05 13:40:35.394010",
    "2021-10-05 13:40:35.394013",
    "2021-10-05 13:40:35.394016",
    "2021-10-05 13:40:35.394020",
    "2021-10-05 13:40:35.394023",
    "2021-10-05 13:40:35.394025",
    "2021-10-05 13:40:35.394029",
    "2021-10-05 13:40:35.394032",
    "2021-10-05 13:40:35.394035",
    "2021-10-05 13:40:35.394038"
)

for ($i = 0; $i -lt $times.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $times[$i]
}
